# MFE_pupae.xlsx update: append new fitness-experiment observations
# (vials at time_hours = 245 and 260) to Sheet1, clear a stray
# font-style override on A115, and leave the selection where the
# author left off (D261) after entering the last data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: row, vial, treatment, time_hours, pupae (pupae omitted -> blank cell)
$rows = @(
    @(205,1,"conditioned",245,5),
    @(206,1,"unconditioned",245,3),
    @(207,2,"conditioned",245,1),
    @(208,2,"unconditioned",245,6),
    @(209,3,"conditioned",245,2),
    @(210,3,"unconditioned",245,5),
    @(211,4,"conditioned",245,4),
    @(212,4,"unconditioned",245,6),
    @(213,5,"conditioned",245,5),
    @(214,5,"unconditioned",245,12),
    @(215,6,"conditioned",245,6),
    @(216,6,"unconditioned",245,7),
    @(217,7,"conditioned",245,2),
    @(218,7,"unconditioned",245,8),
    @(219,8,"conditioned",245,6),
    @(220,8,"unconditioned",245,2),
    @(221,9,"conditioned",245,6),
    @(222,9,"unconditioned",245,4),
    @(223,10,"conditioned",245,1),
    @(224,10,"unconditioned",245,3),
    @(225,11,"conditioned",245,1),
    @(226,11,"unconditioned",245,2),
    @(227,12,"conditioned",245,1),
    @(228,12,"unconditioned",245,2),
    @(229,13,"conditioned",245,1),
    @(230,13,"unconditioned",245,6),
    @(231,14,"conditioned",245,5),
    @(232,14,"unconditioned",245,0),
    @(233,15,"unconditioned",245,3),
    @(234,1,"conditioned",260,6),
    @(235,1,"unconditioned",260,7),
    @(236,2,"conditioned",260,3),
    @(237,2,"unconditioned",260,6),
    @(238,3,"conditioned",260,2),
    @(239,3,"unconditioned",260,2),
    @(240,4,"conditioned",260,4),
    @(241,4,"unconditioned",260,8),
    @(242,5,"conditioned",260,3),
    @(243,5,"unconditioned",260,1),
    @(244,6,"conditioned",260,5),
    @(245,6,"unconditioned",260,11),
    @(246,7,"conditioned",260,9),
    @(247,7,"unconditioned",260,5),
    @(248,8,"conditioned",260,3),
    @(249,8,"unconditioned",260,4),
    @(250,9,"conditioned",260,4),
    @(251,9,"unconditioned",260,0),
    @(252,10,"conditioned",260,6),
    @(253,10,"unconditioned",260,4),
    @(254,11,"conditioned",260,4),
    @(255,11,"unconditioned",260,1),
    @(256,12,"conditioned",260,4),
    @(257,12,"unconditioned",260,0),
    @(258,13,"conditioned",260,4),
    @(259,13,"unconditioned",260,5),
    @(260,14,"conditioned",260,12),
    @(261,14,"unconditioned",260,$null),
    @(262,15,"unconditioned",260,$null)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    if ($row[4] -ne $null) {
        $ws.Cells.Item($r, 4).Value = $row[4]
    }
}

# A115 previously carried an explicit (unused) font-applying style;
# restore it to the sheet's default/normal style.
$ws.Cells.Item(115, 1).Style = "Normal"

# Leave the cursor where the author left it: on the last entered cell.
$ws.Range("D261").Select() | Out-Null
